$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells remain stored as text (matches original inline-string cells)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.172.12"
$ws.Range("E2").Value = "  +6.14%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.916.05"
$ws.Range("E3").Value = "  +2.84%  "
$ws.Range("E4").Value = "  -0.73%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "329.90"
$ws.Range("E5").Value = "  +4.80%  "
$ws.Range("E6").Value = "  -0.75%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5223"
$ws.Range("E7").Value = "  +3.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4083"
$ws.Range("E8").Value = "  +4.73%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08529"
$ws.Range("E9").Value = "  +2.64%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.93"
$ws.Range("E10").Value = "  +1.30%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.124"
$ws.Range("E11").Value = "  +2.14%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.51"
$ws.Range("E12").Value = "  +11.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.456"
$ws.Range("E13").Value = "  +4.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.907.81"
$ws.Range("E14").Value = "  +2.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.403"
$ws.Range("E15").Value = "  +2.55%  "
$ws.Range("E16").Value = "  -0.72%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "95.10"
$ws.Range("E17").Value = "  +4.41%  "
$ws.Range("E18").Value = "  +1.50%  "
$ws.Range("E19").Value = "  -0.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.37"
$ws.Range("E20").Value = "  +4.58%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9997"
$ws.Range("E21").Value = "  -0.75%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.011"
$ws.Range("E22").Value = "  +2.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.195.80"
$ws.Range("E23").Value = "  +6.00%  "
$ws.Range("E24").Value = "  +2.55%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.217"
$ws.Range("E25").Value = "  +0.83%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.135.66"
$ws.Range("E26").Value = "  +2.72%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.09"
$ws.Range("E27").Value = "  +2.15%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.14"
$ws.Range("E28").Value = "  +2.63%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.424"
$ws.Range("E29").Value = "  +0.73%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "128.97"
$ws.Range("E30").Value = "  +2.70%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.086"
$ws.Range("E31").Value = "  +5.24%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.026"
$ws.Range("E33").Value = "  +5.12%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.621"
$ws.Range("E34").Value = "  +0.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02493"
$ws.Range("E35").Value = "  +2.24%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06575"
$ws.Range("E36").Value = "  +0.43%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2207"
$ws.Range("E37").Value = "  +2.67%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.230"
$ws.Range("E38").Value = "  +4.73%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.180"
$ws.Range("E39").Value = "  +3.56%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.879"
$ws.Range("E40").Value = "  -0.54%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6541"
$ws.Range("E41").Value = "  +3.33%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.65"
$ws.Range("E42").Value = "  +5.48%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.244"
$ws.Range("E43").Value = "  +1.25%  "
$ws.Range("E44").Value = "  +2.87%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.23"
$ws.Range("E45").Value = "  +2.18%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.752"
$ws.Range("E46").Value = "  +2.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.084"
$ws.Range("E47").Value = "  +4.89%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.245"
$ws.Range("E48").Value = "  +3.47%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "124.34"
$ws.Range("E49").Value = "  +1.89%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.167"
$ws.Range("E50").Value = "  +2.68%  "
$ws.Range("E51").Value = "  +4.74%  "
